# "création vue initialisation projet"
# Shift the intake-date column (A) for rows 3-63 from the 2015 cohort
# to the 2017 cohort (2015xxxx -> 2017xxxx, i.e. +20000), and reselect
# the working range A1:C63 on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 3; $r -le 63; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value2 = $cell.Value2 + 20000
}

$ws.Range("A1:C63").Select() | Out-Null
